# Update the "view count" style figures in column F for a handful of
# events that appear on both the "展览" (Exhibition) sheet and the
# "全部类型" (All Types) aggregate sheet, reflecting a re-scrape of the
# source site.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# Row -> new F value on the "展览" sheet.
$exhibitionUpdates = @{
    4  = 1281
    9  = 6725
    13 = 6385
    16 = 4300
    19 = 4261
    20 = 218
    32 = 7764
    34 = 1309
    40 = 1545
    42 = 886
    44 = 3854
    49 = 1071
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F value on the "全部类型" sheet (same events, different rows).
$allTypesUpdates = @{
    7  = 1281
    13 = 6725
    17 = 6385
    20 = 4300
    21 = 4261
    22 = 218
    31 = 7764
    33 = 1309
    39 = 1545
    41 = 886
    43 = 3854
    48 = 1071
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
